$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: fill in the new "Use cases" entry under the first effort table
$ws.Range("A9").Value = (Get-Date -Year 2019 -Month 10 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B9").Value = "Use cases"
$ws.Range("C9").Value = 4

# Update selection to reflect the new active cell (cosmetic, matches diff)
$ws.Range("E12").Select()
